$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.723.87"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "3.799.83"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D5").Value = "'706.75"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").Value = "'170.43"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D7").Value = "3.798.74"
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").Value = "'7.38"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "4.441.67"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "3.807.97"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "70.743.22"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'17.39"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "'493.91"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").Value = "'85.06"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").Value = "'12.09"
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("D27").Value = "'10.48"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Value = "3.951.86"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'2.06"
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("D31").Value = "'3.09"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").Value = "'7.33"
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("D34").Value = "'29.11"
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("E35").Value = "  -4.73%  "
$ws.Range("D36").Value = "3.772.19"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "'9.07"
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").Value = "'3.29"
$ws.Range("E43").Value = "  -3.95%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'0.000318"
$ws.Range("E46").Value = "  +3.25%  "
$ws.Range("D47").Value = "'164.35"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "'48.86"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "'423.68"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").Value = "'8.64"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "'1.36"
$ws.Range("E51").Value = "  -1.76%  "
